# Log file updated, with links of Post69
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row 79 below the last data row, copying the formatting of row 78
$ws.Rows("78:78").Copy()
$ws.Rows("79:79").Insert(-4121)  # xlShiftDown

# New row 79 data for Post 69 (Select Loop | Shell Scripting)
$ws.Range("B79").Value = 69
$ws.Range("D79").Value = 44182
$ws.Range("E79").Value = "https://programmingport.hashnode.dev/select-loop-or-shell-scripting"
$ws.Range("C79").Value = "Select Loop | Shell Scripting"
$ws.Range("F79").Value = "https://dev.to/rahulmishra05/select-loop-shell-scripting-3lme"

# Extend the table range to include the new row
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("B10:F79"))

$ws.Range("F79").Select()
